$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 'CasesTab' query text (Cohort column removed from the Cypher query/output).
# Re-assigning this text causes the shared-string table to be rebuilt: the old
# (now-unused) Cases-with-Cohort string is dropped, the Sample/Files query strings
# shift up, and this new string is appended last - matching the target ordering.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Border Collie','Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder', 'Bladder, Urethra']
    
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $casesQuery

# Row heights were re-autofit by Excel after the text changed (wrap text is on for B:C).
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 275.5

# Selection / scroll position moved to B2 on save.
$ws.Range("B2").Select() | Out-Null

